$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Formatting first (copy/paste of formats does not touch the shared
# string table, so this can happen before any of the text is typed in).
# ------------------------------------------------------------------

# Rows 14-16: Date/No columns get the same formats as the row above them.
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A14:A16").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14:B16").PasteSpecial(-4122) | Out-Null

# D14 and D16 hold long, multi-line notes, so they need the wrap-text
# style already used by the other multi-line Detail cells (D6/D12).
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Text values, entered in the same order the shared-string table ends
# up in.
# ------------------------------------------------------------------
$ws.Range("D3").Value = "npm i"
$ws.Range("D5").Value = "git add ."
$ws.Range("D7").Value = "git commit -m `"First Commit for project red-center-project`""

$ws.Range("D14").Value = "git add .`n.gitignore`nprojectHistory.xlsx`nsrc/app/@theme/components/header/header.component.html`nsrc/app/@theme/components/header/header.component.ts`nsrc/index.html"

$ws.Range("C16").Value = "Merge 'feature/header' into 'develop' branch "

$ws.Range("C9").Value = "Switch to the new branch by using 'git switch -c `"feature/header`"'"
$ws.Range("D9").Value = "git switch -c `"feature/header`""

$ws.Range("D16").Value = "git checkout develop`ngit merge feature/header"

$ws.Range("C15").Value = "Commit files into GIT by using 'git commit -m `"Revise Title on index.html and header.component feature with update projectHistory.xlsx`"'"
$ws.Range("D15").Value = "git commit -m `"Revise Title on index.html and header.component feature with update projectHistory.xlsx`""

# Remaining cells reuse existing shared strings / are plain numbers, so
# their order has no effect on the shared-string table layout.
$ws.Range("C14").Value = "Add files by using 'git add .'"

$ws.Range("A14").Value = 44986
$ws.Range("B14").Value = 13
$ws.Range("A15").Value = 44986
$ws.Range("B15").Value = 14
$ws.Range("A16").Value = 44986
$ws.Range("B16").Value = 15

$ws.Rows(14).RowHeight = 105
$ws.Rows(16).RowHeight = 30

# ------------------------------------------------------------------
# Grow the table to cover the three new rows and refresh the selection.
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D16"))

$ws.Range("D16").Select() | Out-Null
